$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3410236666666666
$ws.Range("H2").Value = 1.023071
$ws.Range("I2").Value = 0.01850325494520333
$ws.Range("J2").Value = 0.01850325494520333
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04191666666666666
$ws.Range("N2").Value = 0.12575
$ws.Range("O2").Value = 0.007985764192544619
$ws.Range("P2").Value = 0.007985764192544619
$ws.Range("Q2").Value = 0.01429457536111111
$ws.Range("R2").Value = 0.12865117825
$ws.Range("S2").Value = 0.0001477626307869289
$ws.Range("T2").Value = 0.0001477626307869289

$ws.Range("G3").Value = 0.3410236666666666
$ws.Range("H3").Value = 1.023071
$ws.Range("I3").Value = 0.01850325494520333
$ws.Range("J3").Value = 0.01850325494520333
$ws.Range("M3").Value = 5.207007
$ws.Range("N3").Value = 15.621021
$ws.Range("O3").Value = 0.9920142358074554
$ws.Range("P3").Value = 0.9920142358074554
$ws.Range("Q3").Value = 1.775712619499
$ws.Range("R3").Value = 15.981413575491
$ws.Range("S3").Value = 0.0183554923144164
$ws.Range("T3").Value = 0.0183554923144164

$ws.Range("I4").Value = 0.2085050756621187
$ws.Range("J4").Value = 0.2085050756621187
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04191666666666666
$ws.Range("N4").Value = 0.12575
$ws.Range("O4").Value = 0.007985764192544619
$ws.Range("P4").Value = 0.007985764192544619
$ws.Range("Q4").Value = 0.1610793088055556
$ws.Range("R4").Value = 1.44971377925
$ws.Range("S4").Value = 0.001665072367186354
$ws.Range("T4").Value = 0.001665072367186354

$ws.Range("I5").Value = 0.2085050756621187
$ws.Range("J5").Value = 0.2085050756621187
$ws.Range("M5").Value = 5.207007
$ws.Range("N5").Value = 15.621021
$ws.Range("O5").Value = 0.9920142358074554
$ws.Range("P5").Value = 0.9920142358074554
$ws.Range("Q5").Value = 20.009727757591
$ws.Range("R5").Value = 180.087549818319
$ws.Range("S5").Value = 0.2068400032949324
$ws.Range("T5").Value = 0.2068400032949323

$ws.Range("G6").Value = 0.2092423333333333
$ws.Range("H6").Value = 0.627727
$ws.Range("I6").Value = 0.0113530661283407
$ws.Range("J6").Value = 0.0113530661283407
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04191666666666666
$ws.Range("N6").Value = 0.12575
$ws.Range("O6").Value = 0.007985764192544619
$ws.Range("P6").Value = 0.007985764192544619
$ws.Range("Q6").Value = 0.008770741138888889
$ws.Range("R6").Value = 0.07893667025000001
$ws.Range("S6").Value = 0.00009066290896329434
$ws.Range("T6").Value = 0.00009066290896329434

$ws.Range("G7").Value = 0.2092423333333333
$ws.Range("H7").Value = 0.627727
$ws.Range("I7").Value = 0.0113530661283407
$ws.Range("J7").Value = 0.0113530661283407
$ws.Range("M7").Value = 5.207007
$ws.Range("N7").Value = 15.621021
$ws.Range("O7").Value = 0.9920142358074554
$ws.Range("P7").Value = 0.9920142358074554
$ws.Range("Q7").Value = 1.089526294363
$ws.Range("R7").Value = 9.805736649267
$ws.Range("S7").Value = 0.01126240321937741
$ws.Range("T7").Value = 0.01126240321937741

$ws.Range("G8").Value = 14.03735666666667
$ws.Range("H8").Value = 42.11207
$ws.Range("I8").Value = 0.7616386032643372
$ws.Range("J8").Value = 0.7616386032643372
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04191666666666666
$ws.Range("N8").Value = 0.12575
$ws.Range("O8").Value = 0.007985764192544619
$ws.Range("P8").Value = 0.007985764192544619
$ws.Range("Q8").Value = 0.5883992002777778
$ws.Range("R8").Value = 5.295592802500001
$ws.Range("S8").Value = 0.006082266285608041
$ws.Range("T8").Value = 0.006082266285608041

$ws.Range("G9").Value = 14.03735666666667
$ws.Range("H9").Value = 42.11207
$ws.Range("I9").Value = 0.7616386032643372
$ws.Range("J9").Value = 0.7616386032643372
$ws.Range("M9").Value = 5.207007
$ws.Range("N9").Value = 15.621021
$ws.Range("O9").Value = 0.9920142358074554
$ws.Range("P9").Value = 0.9920142358074554
$ws.Range("Q9").Value = 73.09261442483
$ws.Range("R9").Value = 657.83352982347
$ws.Range("S9").Value = 0.7555563369787291
$ws.Range("T9").Value = 0.7555563369787291

